$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style s="1") from the last existing data row (A205) down through the new rows
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new year index (column A) and normalized value (column B) data
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = [double]"9.251858538542971E-18"
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = [double]"1.009293658750142E-17"
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = [double]"-7.494005416219807E-17"
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = [double]"-1.233581138472396E-17"
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = [double]"-9.71445146547012E-17"
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = [double]"4.361590453884543E-17"
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = [double]"-2.544261098099317E-17"
$ws.Range("A213").Value = 211
$ws.Range("B213").Value = [double]"9.159339953157541E-17"
$ws.Range("A214").Value = 212
$ws.Range("B214").Value = [double]"1.110223024625157E-16"
$ws.Range("A215").Value = 213
$ws.Range("B215").Value = [double]"0"
$ws.Range("A216").Value = 214
$ws.Range("B216").Value = [double]"0"
$ws.Range("A217").Value = 215
$ws.Range("B217").Value = [double]"0"
